$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "nbk-vl"
$ws.Cells.Item(2, 2).Value = "Truong trung hoc Chuyen Nguyen Binh Khiem"
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(2, 4).Value = "Vinh Long"

$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = "337829999"
$ws.Cells.Item(2, 5).ClearFormats()

# Row 3
$ws.Cells.Item(3, 1).Value = "nbk-qn"
$ws.Cells.Item(3, 2).Value = "Truong trung hoc Chuyen Nguyen Binh Khiem"
$ws.Cells.Item(3, 3).Value = 1
$ws.Cells.Item(3, 4).Value = "Quang Ngai"

$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = "334442222"
$ws.Cells.Item(3, 5).ClearFormats()

# Column width adjustments (auto-fit after new, wider data was added).
# (Values chosen land the engine's internal 1/6-character rounding on the
# stored width closest to the canonical bestFit widths.)
$ws.Columns.Item(1).ColumnWidth = 6.5
$ws.Columns.Item(2).ColumnWidth = 40.833333333333336
$ws.Columns.Item(4).ColumnWidth = 10.666666666666666
